$wb = $excel.ActiveWorkbook

# Sheet1 "nhap-linhkien": add a new row of import data
$ws1 = $wb.Worksheets.Item("nhap-linhkien")
$ws1.Range("A2").Value = "GRM219R61A106KE44D"
$ws1.Range("B2").Value = "Ceramic Cap 0805 10uF 10V 10% X5R"
$ws1.Range("C2").Value = "sohopdong01"
$ws1.Range("D2").Value = "sanpham01"
$ws1.Range("E2").Value = "cty01"
$ws1.Range("F2").NumberFormat = "@"
$ws1.Range("F2").Value = "2021-09-11"
$ws1.Range("G2").Value = "Cái"
$ws1.Range("H2").Value = 12
$ws1.Range("I2").Value = 21
$ws1.Range("J2").Value = 252

# Sheet3 "ton-linhkien": add corresponding stock row
$ws3 = $wb.Worksheets.Item("ton-linhkien")
$ws3.Range("A2").Value = "Ceramic Cap 0805 10uF 10V 10% X5R"
$ws3.Range("B2").Value = 12
$ws3.Range("C2").Value = "Cái"
